# Insert a new data row at row 721 (2026/01/25, 日, 13:00, rank 161),
# pushing the existing rows 721-762 (2026/12/29 .. 2027/01/05) down to
# 722-763. This mirrors the XML diff: <dimension ref="A1:D762"/> ->
# <dimension ref="A1:D763"/> plus the shifted row data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 721.. down by one to make room for the new record.
$ws.Rows("721").Insert()

$dateCell = $ws.Range("A721")

# Column A stores dates as plain text (e.g. "2026/12/29"), not real
# Excel date serials. Assigning a date-shaped string normally gets
# auto-converted to a date value/format, so force a Text number format
# first, then clear the formatting afterwards so the cell ends up back
# at the default style (matching every other data row) while keeping
# the literal text content.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/01/25"
$dateCell.ClearFormats()

$ws.Range("B721").Value = "日"
$ws.Range("C721").Value = 13
$ws.Range("D721").Value = 161
